$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.174.73"
$ws.Range("E2").Value = "  -0.07%  "

$ws.Range("D3").Value = "1.866.59"
$ws.Range("E3").Value = "  +0.40%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.41"
$ws.Range("E5").Value = "  -0.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4690"
$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2854"
$ws.Range("E8").Value = "  -1.45%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "41.49"
$ws.Range("E9").Value = "  -2.86%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06558"
$ws.Range("E10").Value = "  -0.41%  "

$ws.Range("E11").Value = "  -2.23%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07770"

$ws.Range("D13").Value = "1.873.71"
$ws.Range("E13").Value = "  +0.82%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "95.85"
$ws.Range("E14").Value = "  -1.80%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6915"
$ws.Range("E15").Value = "  +1.91%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.102"
$ws.Range("E16").Value = "  -0.28%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "265.87"
$ws.Range("E17").Value = "  -1.20%  "

$ws.Range("D18").Value = "30.168.05"
$ws.Range("E18").Value = "  -0.02%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.65"
$ws.Range("E19").Value = "  +0.20%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007727"
$ws.Range("E20").Value = "  +0.74%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9997"
$ws.Range("E21").Value = "  +0.08%  "

$ws.Range("D22").Value = "2.125.86"
$ws.Range("E22").Value = "  +1.31%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9993"
$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.258"
$ws.Range("E24").Value = "  +0.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.159"
$ws.Range("E25").Value = "  -0.29%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.469"
$ws.Range("E26").Value = "  +2.94%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.58"
$ws.Range("E27").Value = "  -0.35%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.67"
$ws.Range("E28").Value = "  -1.56%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.936"
$ws.Range("E29").Value = "  -0.53%  "

$ws.Range("E30").Value = "  -0.44%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09928"
$ws.Range("E31").Value = "  +0.22%  "

$ws.Range("E32").Value = "  +0.55%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.460"
$ws.Range("E33").Value = "  -0.48%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.049"
$ws.Range("E34").Value = "  +0.23%  "

$ws.Range("E35").Value = "  +0.25%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.128"
$ws.Range("E36").Value = "  +0.45%  "

$ws.Range("E37").Value = "  -0.30%  "

$ws.Range("E38").Value = "  +0.60%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01856"
$ws.Range("E39").Value = "  -1.06%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.774"
$ws.Range("E40").Value = "  +6.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.305"
$ws.Range("E41").Value = "  -0.60%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.32"
$ws.Range("E42").Value = "  -1.70%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.931"
$ws.Range("E43").Value = "  -0.32%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9995"
$ws.Range("E44").Value = "  +0.12%  "

$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4141"
$ws.Range("E45").Value = "  -0.31%  "

$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.8359"
$ws.Range("E46").Value = "  -0.45%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.98"
$ws.Range("E47").Value = "  -0.64%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "970.00"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.110"
$ws.Range("E49").Value = "  +0.24%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.127"
$ws.Range("E50").Value = "  -1.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.60"
$ws.Range("E51").Value = "  +1.52%  "
